$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename the two updated sheets (new end-date of 12/14/22 instead of 9/25/22)
# ---------------------------------------------------------------------------
$wsJulySept = $wb.Worksheets.Item("July 21 - Sept 22")
$wsJulySept.Name = "July 21 - Dec 22"

$wsAprSept = $wb.Worksheets.Item("April 21 - Sept 22")
$wsAprSept.Name = "April 21 - Dec 22"

$wsJulyDec = $wb.Worksheets.Item("July 21 - Dec 22")
$wsAprDec = $wb.Worksheets.Item("April 21 - Dec 22")

# ---------------------------------------------------------------------------
# 2. "July 21 - Dec 22" sheet: add the newest data point (new row 10), which
#    shifts the totals row and everything below it down by one row.
# ---------------------------------------------------------------------------

# Insert a new row at 10 (pushes old row 10 -> 11, etc.)
$wsJulyDec.Rows.Item(10).Insert()

# Copy the formatting from row 9 (the prior data row) into the new row 10
$wsJulyDec.Range("A9:F9").Copy()
$wsJulyDec.Range("A10:F10").PasteSpecial(-4122)
$wsJulyDec.Rows.Item(10).RowHeight = $wsJulyDec.Rows.Item(9).RowHeight

# New data point: date 10/23/2022 (serial 44857)
$wsJulyDec.Range("A10").Value = 44857
$wsJulyDec.Range("B10").Value = 1.17
$wsJulyDec.Range("C10").Value = 0.04
$wsJulyDec.Range("D10").Value = 364
$wsJulyDec.Range("E10").Formula = "=D10/D11"
$wsJulyDec.Range("F10").Formula = "=C10*E10"

# Totals row (now row 11) needs its SUM ranges extended to include row 10
$wsJulyDec.Range("D11").Formula = "=SUM(D3:D10)"
$wsJulyDec.Range("E11").Formula = "=SUM(E3:E10)"
$wsJulyDec.Range("F11").Formula = "=SUM(F3:F10)"

# Updated "Deaths at end of period" figure (now row 17 after the row insert)
$wsJulyDec.Range("B17").Value = 1083279

# Title text update (row 1)
$wsJulyDec.Range("A1").Value = "Preventable COVID-19 Deaths, with available vaccination, 7/15/21 to 12/14/22"

# Booster definition note text update (now row 31 after the row insert)
$wsJulyDec.Range("B31").Value = "Vaccinated = primary series + all boosters recommended."

# ---------------------------------------------------------------------------
# 3. "April 21 - Dec 22" sheet: update deaths figure + dependent formulas/text
# ---------------------------------------------------------------------------
$wsAprDec.Range("B4").Value = 1083279
$wsAprDec.Range("A1").Value = "Preventable COVID-19 Deaths, with available vaccination, 4/15/21 to 12/14/22"
$wsAprDec.Range("A8").Value = "Preventable deaths 7/15/21 to 12/14/22, estimated"
$wsAprDec.Range("B8").Formula = "='July 21 - Dec 22'!B20"

Write-Output "done"
